$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I (I0) and J (IF), copying the formatting
# of the existing header cell H1 (bold, bordered, centered style).
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# I0/IF data values for rows 2-75 (columns I and J respectively).
$data = @(
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 10),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(10, 10),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(5, 5),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(4, 4),
    @(5, 5)
)

for ($k = 0; $k -lt $data.Count; $k++) {
    $row = $k + 2
    $ws.Cells.Item($row, 9).Value = $data[$k][0]
    $ws.Cells.Item($row, 10).Value = $data[$k][1]
}

Write-Output "I0 and IF columns added"
